$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 42

# New buy row appended below the existing data (08/20/2025 run).
# Column A holds the date as literal text (matches the other recent rows,
# which are also plain "MM/DD/YYYY" strings rather than Excel date serials).
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "08/20/2025"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 0.0004368600000000007
$ws.Cells.Item($row, 3).Value = 114453.1428833034
$ws.Cells.Item($row, 4).Value = 50
